$d = $word.ActiveDocument

function Merge-IdRun($idNum) {
    $old = "<id>p126v_$idNum</id>"
    $new = "<id>p126v_$idNum</id>"
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Merge-IdRun 1
Merge-IdRun 2
Merge-IdRun 3
